$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Word-code"

$ws.Cells.Item(2, 2).Value = "def"
$ws.Cells.Item(2, 3).Value = 8
$ws.Cells.Item(3, 2).Value = "getAllWordCode"
$ws.Cells.Item(3, 3).Value = 81
$ws.Cells.Item(4, 2).Value = "("
$ws.Cells.Item(4, 3).Value = 61
$ws.Cells.Item(5, 2).Value = "self"
$ws.Cells.Item(5, 3).Value = 5
$ws.Cells.Item(6, 2).Value = ","
$ws.Cells.Item(6, 3).Value = 34
$ws.Cells.Item(7, 2).Value = "data_list"
$ws.Cells.Item(7, 3).Value = 81
$ws.Cells.Item(8, 2).Value = ")"
$ws.Cells.Item(8, 3).Value = 62
$ws.Cells.Item(9, 2).Value = ":"
$ws.Cells.Item(9, 3).Value = 33
$ws.Cells.Item(10, 2).Value = "word_code_list"
$ws.Cells.Item(10, 3).Value = 81
$ws.Cells.Item(11, 2).Value = "'="
$ws.Cells.Item(11, 2).Style = "Normal"
$ws.Cells.Item(11, 3).Value = 46
$ws.Cells.Item(12, 2).Value = "["
$ws.Cells.Item(12, 3).Value = 63
$ws.Cells.Item(13, 2).Value = "]"
$ws.Cells.Item(13, 3).Value = 64
$ws.Cells.Item(14, 2).Value = "for"
$ws.Cells.Item(14, 3).Value = 13
$ws.Cells.Item(15, 2).Value = "sentence"
$ws.Cells.Item(15, 3).Value = 81
$ws.Cells.Item(16, 2).Value = "in"
$ws.Cells.Item(16, 3).Value = 72
$ws.Cells.Item(17, 2).Value = "data_list"
$ws.Cells.Item(17, 3).Value = 81
$ws.Cells.Item(18, 2).Value = ":"
$ws.Cells.Item(18, 3).Value = 33
$ws.Cells.Item(19, 2).Value = "#"
$ws.Cells.Item(19, 3).Value = 43
$ws.Cells.Item(20, 2).Value = "print"
$ws.Cells.Item(20, 3).Value = 25
$ws.Cells.Item(21, 2).Value = "("
$ws.Cells.Item(21, 3).Value = 61
$ws.Cells.Item(22, 2).Value = "sentence"
$ws.Cells.Item(22, 3).Value = 81
$ws.Cells.Item(23, 2).Value = ")"
$ws.Cells.Item(23, 3).Value = 62
$ws.Cells.Item(24, 2).Value = "sen_list"
$ws.Cells.Item(24, 3).Value = 81
$ws.Cells.Item(25, 2).Value = "'="
$ws.Cells.Item(25, 2).Style = "Normal"
$ws.Cells.Item(25, 3).Value = 46
$ws.Cells.Item(26, 2).Value = "sentence"
$ws.Cells.Item(26, 3).Value = 81
$ws.Cells.Item(27, 2).Value = "."
$ws.Cells.Item(27, 3).Value = 32
$ws.Cells.Item(28, 2).Value = "split"
$ws.Cells.Item(28, 3).Value = 81
$ws.Cells.Item(29, 2).Value = "("
$ws.Cells.Item(29, 3).Value = 61
$ws.Cells.Item(30, 2).Value = "''"
$ws.Cells.Item(30, 2).Style = "Normal"
$ws.Cells.Item(30, 3).Value = 67
$ws.Cells.Item(31, 2).Value = "''"
$ws.Cells.Item(31, 2).Style = "Normal"
$ws.Cells.Item(31, 3).Value = 67
$ws.Cells.Item(32, 2).Value = ")"
$ws.Cells.Item(32, 3).Value = 62
$ws.Cells.Item(33, 2).Value = "#"
$ws.Cells.Item(33, 3).Value = 43
$ws.Cells.Item(34, 2).Value = "print"
$ws.Cells.Item(34, 3).Value = 25
$ws.Cells.Item(35, 2).Value = "("
$ws.Cells.Item(35, 3).Value = 61
$ws.Cells.Item(36, 2).Value = "sen_list"
$ws.Cells.Item(36, 3).Value = 81
$ws.Cells.Item(37, 2).Value = ")"
$ws.Cells.Item(37, 3).Value = 62
$ws.Cells.Item(38, 2).Value = "#"
$ws.Cells.Item(38, 3).Value = 43
$ws.Cells.Item(39, 2).Value = "print"
$ws.Cells.Item(39, 3).Value = 25
$ws.Cells.Item(40, 2).Value = "("
$ws.Cells.Item(40, 3).Value = 61
$ws.Cells.Item(41, 2).Value = "sen_list"
$ws.Cells.Item(41, 3).Value = 81
$ws.Cells.Item(42, 2).Value = ")"
$ws.Cells.Item(42, 3).Value = 62
$ws.Cells.Item(43, 2).Value = "tmp_list_origin"
$ws.Cells.Item(43, 3).Value = 81
$ws.Cells.Item(44, 2).Value = "'="
$ws.Cells.Item(44, 2).Style = "Normal"
$ws.Cells.Item(44, 3).Value = 46
$ws.Cells.Item(45, 2).Value = "["
$ws.Cells.Item(45, 3).Value = 63
$ws.Cells.Item(46, 2).Value = "]"
$ws.Cells.Item(46, 3).Value = 64
$ws.Cells.Item(47, 2).Value = "tmp_list_dict"
$ws.Cells.Item(47, 3).Value = 81
$ws.Cells.Item(48, 2).Value = "'="
$ws.Cells.Item(48, 2).Style = "Normal"
$ws.Cells.Item(48, 3).Value = 46
$ws.Cells.Item(49, 2).Value = "["
$ws.Cells.Item(49, 3).Value = 63
$ws.Cells.Item(50, 2).Value = "]"
$ws.Cells.Item(50, 3).Value = 64
$ws.Cells.Item(51, 2).Value = "for"
$ws.Cells.Item(51, 3).Value = 13
$ws.Cells.Item(52, 2).Value = "word"
$ws.Cells.Item(52, 3).Value = 81
$ws.Cells.Item(53, 2).Value = "in"
$ws.Cells.Item(53, 3).Value = 72
$ws.Cells.Item(54, 2).Value = "sen_list"
$ws.Cells.Item(54, 3).Value = 81
$ws.Cells.Item(55, 2).Value = ":"
$ws.Cells.Item(55, 3).Value = 33
$ws.Cells.Item(56, 2).Value = "#"
$ws.Cells.Item(56, 3).Value = 43
$ws.Cells.Item(57, 2).Value = "print"
$ws.Cells.Item(57, 3).Value = 25
$ws.Cells.Item(58, 2).Value = "("
$ws.Cells.Item(58, 3).Value = 61
$ws.Cells.Item(59, 2).Value = "word"
$ws.Cells.Item(59, 3).Value = 81
$ws.Cells.Item(60, 2).Value = ")"
$ws.Cells.Item(60, 3).Value = 62
$ws.Cells.Item(61, 2).Value = "tmp_list_origin"
$ws.Cells.Item(61, 3).Value = 81
$ws.Cells.Item(62, 2).Value = "."
$ws.Cells.Item(62, 3).Value = 32
$ws.Cells.Item(63, 2).Value = "extend"
$ws.Cells.Item(63, 3).Value = 81
$ws.Cells.Item(64, 2).Value = "("
$ws.Cells.Item(64, 3).Value = 61
$ws.Cells.Item(65, 2).Value = "self"
$ws.Cells.Item(65, 3).Value = 5
$ws.Cells.Item(66, 2).Value = "."
$ws.Cells.Item(66, 3).Value = 32
$ws.Cells.Item(67, 2).Value = "__stringProcessing"
$ws.Cells.Item(67, 3).Value = 81
$ws.Cells.Item(68, 2).Value = "("
$ws.Cells.Item(68, 3).Value = 61
$ws.Cells.Item(69, 2).Value = "word"
$ws.Cells.Item(69, 3).Value = 81
$ws.Cells.Item(70, 2).Value = ")"
$ws.Cells.Item(70, 3).Value = 62
$ws.Cells.Item(71, 2).Value = ")"
$ws.Cells.Item(71, 3).Value = 62
$ws.Cells.Item(72, 2).Value = "#"
$ws.Cells.Item(72, 3).Value = 43
$ws.Cells.Item(73, 2).Value = "print"
$ws.Cells.Item(73, 3).Value = 25
$ws.Cells.Item(74, 2).Value = "("
$ws.Cells.Item(74, 3).Value = 61
$ws.Cells.Item(75, 2).Value = "tmp_list_origin"
$ws.Cells.Item(75, 3).Value = 81
$ws.Cells.Item(76, 2).Value = ")"
$ws.Cells.Item(76, 3).Value = 62
$ws.Cells.Item(77, 2).Value = "for"
$ws.Cells.Item(77, 3).Value = 13
$ws.Cells.Item(78, 2).Value = "word"
$ws.Cells.Item(78, 3).Value = 81
$ws.Cells.Item(79, 2).Value = "in"
$ws.Cells.Item(79, 3).Value = 72
$ws.Cells.Item(80, 2).Value = "tmp_list_origin"
$ws.Cells.Item(80, 3).Value = 81
$ws.Cells.Item(81, 2).Value = ":"
$ws.Cells.Item(81, 3).Value = 33
$ws.Cells.Item(82, 2).Value = "tmp_list_dict"
$ws.Cells.Item(82, 3).Value = 81
$ws.Cells.Item(83, 2).Value = "."
$ws.Cells.Item(83, 3).Value = 32
$ws.Cells.Item(84, 2).Value = "append"
$ws.Cells.Item(84, 3).Value = 81
$ws.Cells.Item(85, 2).Value = "("
$ws.Cells.Item(85, 3).Value = 61
$ws.Cells.Item(86, 2).Value = "{"
$ws.Cells.Item(86, 3).Value = 65
$ws.Cells.Item(87, 2).Value = "word"
$ws.Cells.Item(87, 3).Value = 81
$ws.Cells.Item(88, 2).Value = ":"
$ws.Cells.Item(88, 3).Value = 33
$ws.Cells.Item(89, 2).Value = "self"
$ws.Cells.Item(89, 3).Value = 5
$ws.Cells.Item(90, 2).Value = "."
$ws.Cells.Item(90, 3).Value = 32
$ws.Cells.Item(91, 2).Value = "__getWordCode"
$ws.Cells.Item(91, 3).Value = 81
$ws.Cells.Item(92, 2).Value = "("
$ws.Cells.Item(92, 3).Value = 61
$ws.Cells.Item(93, 2).Value = "word"
$ws.Cells.Item(93, 3).Value = 81
$ws.Cells.Item(94, 2).Value = ")"
$ws.Cells.Item(94, 3).Value = 62
$ws.Cells.Item(95, 2).Value = "}"
$ws.Cells.Item(95, 3).Value = 66
$ws.Cells.Item(96, 2).Value = ")"
$ws.Cells.Item(96, 3).Value = 62
$ws.Cells.Item(97, 2).Value = "word_code_list"
$ws.Cells.Item(97, 3).Value = 81
$ws.Cells.Item(98, 2).Value = "."
$ws.Cells.Item(98, 3).Value = 32
$ws.Cells.Item(99, 2).Value = "extend"
$ws.Cells.Item(99, 3).Value = 81
$ws.Cells.Item(100, 2).Value = "("
$ws.Cells.Item(100, 3).Value = 61
$ws.Cells.Item(101, 2).Value = "tmp_list_dict"
$ws.Cells.Item(101, 3).Value = 81
$ws.Cells.Item(102, 2).Value = ")"
$ws.Cells.Item(102, 3).Value = 62
$ws.Cells.Item(103, 2).Value = "self"
$ws.Cells.Item(103, 3).Value = 5
$ws.Cells.Item(104, 2).Value = "."
$ws.Cells.Item(104, 3).Value = 32
$ws.Cells.Item(105, 2).Value = "__storeResultsInExcel"
$ws.Cells.Item(105, 3).Value = 81
$ws.Cells.Item(106, 2).Value = "("
$ws.Cells.Item(106, 3).Value = 61
$ws.Cells.Item(107, 2).Value = "word_code_list"
$ws.Cells.Item(107, 3).Value = 81
$ws.Cells.Item(108, 2).Value = ")"
$ws.Cells.Item(108, 3).Value = 62
$ws.Cells.Item(109, 2).Value = "return"
$ws.Cells.Item(109, 3).Value = 27
$ws.Cells.Item(110, 2).Value = "word_code_list"
$ws.Cells.Item(110, 3).Value = 81

$ws.Rows.Item(111).Delete()
